$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 13:05:50"

# zh-cn sheet: "Correspond Handoff Datetime" (col H) and "Correspond Handback DateTime" (col K) for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 13:05:45"
$wsZhCn.Range("K2").Value = "2016-08-26 13:06:08"

# de-de sheet: "Correspond Handoff Datetime" (col H, shares the same underlying
# text as Overview's G2 in the original workbook) and "Correspond Handback
# DateTime" (col K) for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 13:05:50"
$wsDeDe.Range("K2").Value = "2016-08-26 13:06:19"
